$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 15.281766
$ws.Range("H2").Value = 45.845298
$ws.Range("I2").Value = 0.1817381432449346
$ws.Range("J2").Value = 0.1817381432449346
$ws.Range("M2").Value = 46.63275166666667
$ws.Range("N2").Value = 139.898255
$ws.Range("O2").Value = 0.9158911059585902
$ws.Range("P2").Value = 0.9158911059585902
$ws.Range("Q2").Value = 712.63079890611
$ws.Range("R2").Value = 6413.67719015499
$ws.Range("S2").Value = 0.1664523490114639
$ws.Range("T2").Value = 0.1664523490114639

$ws.Range("G3").Value = 15.281766
$ws.Range("H3").Value = 45.845298
$ws.Range("I3").Value = 0.1817381432449346
$ws.Range("J3").Value = 0.1817381432449346
$ws.Range("O3").Value = 0.05441917700612491
$ws.Range("P3").Value = 0.05441917700612491
$ws.Range("Q3").Value = 42.34213143177
$ws.Range("R3").Value = 381.07918288593
$ws.Range("S3").Value = 0.009890040186010579
$ws.Range("T3").Value = 0.009890040186010581

$ws.Range("G4").Value = 15.281766
$ws.Range("H4").Value = 45.845298
$ws.Range("I4").Value = 0.1817381432449346
$ws.Range("J4").Value = 0.1817381432449346
$ws.Range("M4").Value = 0.8496050000000001
$ws.Range("N4").Value = 2.548815
$ws.Range("O4").Value = 0.01668667696558362
$ws.Range("P4").Value = 0.01668667696558362
$ws.Range("Q4").Value = 12.98346480243
$ws.Range("R4").Value = 116.85118322187
$ws.Range("S4").Value = 0.003032605688653188
$ws.Range("T4").Value = 0.003032605688653188

$ws.Range("G5").Value = 15.281766
$ws.Range("H5").Value = 45.845298
$ws.Range("I5").Value = 0.1817381432449346
$ws.Range("J5").Value = 0.1817381432449346
$ws.Range("M5").Value = 0.662052
$ws.Range("N5").Value = 1.986156
$ws.Range("O5").Value = 0.01300304006970129
$ws.Range("P5").Value = 0.0130030400697013
$ws.Range("Q5").Value = 10.117323743832
$ws.Range("R5").Value = 91.055913694488
$ws.Range("S5").Value = 0.002363148358806998
$ws.Range("T5").Value = 0.002363148358806999

$ws.Range("I6").Value = 0.2947137116012682
$ws.Range("J6").Value = 0.2947137116012682
$ws.Range("M6").Value = 46.63275166666667
$ws.Range("N6").Value = 139.898255
$ws.Range("O6").Value = 0.9158911059585902
$ws.Range("P6").Value = 0.9158911059585902
$ws.Range("Q6").Value = 1155.63009502052
$ws.Range("R6").Value = 10400.67085518468
$ws.Range("S6").Value = 0.2699256672596466
$ws.Range("T6").Value = 0.2699256672596466

$ws.Range("I7").Value = 0.2947137116012682
$ws.Range("J7").Value = 0.2947137116012682
$ws.Range("O7").Value = 0.05441917700612491
$ws.Range("P7").Value = 0.05441917700612491
$ws.Range("S7").Value = 0.01603807763776146
$ws.Range("T7").Value = 0.01603807763776147

$ws.Range("I8").Value = 0.2947137116012682
$ws.Range("J8").Value = 0.2947137116012682
$ws.Range("M8").Value = 0.8496050000000001
$ws.Range("N8").Value = 2.548815
$ws.Range("O8").Value = 0.01668667696558362
$ws.Range("P8").Value = 0.01668667696558362
$ws.Range("Q8").Value = 21.05449650276
$ws.Range("R8").Value = 189.49046852484
$ws.Range("S8").Value = 0.004917792502818538
$ws.Range("T8").Value = 0.004917792502818538

$ws.Range("I9").Value = 0.2947137116012682
$ws.Range("J9").Value = 0.2947137116012682
$ws.Range("M9").Value = 0.662052
$ws.Range("N9").Value = 1.986156
$ws.Range("O9").Value = 0.01300304006970129
$ws.Range("P9").Value = 0.0130030400697013
$ws.Range("Q9").Value = 16.406649582624
$ws.Range("R9").Value = 147.659846243616
$ws.Range("S9").Value = 0.003832174201041682
$ws.Range("T9").Value = 0.003832174201041683

$ws.Range("G10").Value = 18.371237
$ws.Range("H10").Value = 55.113711
$ws.Range("I10").Value = 0.2184796247693259
$ws.Range("J10").Value = 0.2184796247693259
$ws.Range("M10").Value = 46.63275166666667
$ws.Range("N10").Value = 139.898255
$ws.Range("O10").Value = 0.9158911059585902
$ws.Range("P10").Value = 0.9158911059585902
$ws.Range("Q10").Value = 856.7013328304785
$ws.Range("R10").Value = 7710.311995474306
$ws.Range("S10").Value = 0.2001035451593957
$ws.Range("T10").Value = 0.2001035451593957

$ws.Range("G11").Value = 18.371237
$ws.Range("H11").Value = 55.113711
$ws.Range("I11").Value = 0.2184796247693259
$ws.Range("J11").Value = 0.2184796247693259
$ws.Range("O11").Value = 0.05441917700612491
$ws.Range("P11").Value = 0.05441917700612491
$ws.Range("Q11").Value = 50.90231924884834
$ws.Range("R11").Value = 458.120873239635
$ws.Range("S11").Value = 0.0118894813725537
$ws.Range("T11").Value = 0.0118894813725537

$ws.Range("G12").Value = 18.371237
$ws.Range("H12").Value = 55.113711
$ws.Range("I12").Value = 0.2184796247693259
$ws.Range("J12").Value = 0.2184796247693259
$ws.Range("M12").Value = 0.8496050000000001
$ws.Range("N12").Value = 2.548815
$ws.Range("O12").Value = 0.01668667696558362
$ws.Range("P12").Value = 0.01668667696558362
$ws.Range("Q12").Value = 15.608294811385
$ws.Range("R12").Value = 140.474653302465
$ws.Range("S12").Value = 0.003645698922087665
$ws.Range("T12").Value = 0.003645698922087665

$ws.Range("G13").Value = 18.371237
$ws.Range("H13").Value = 55.113711
$ws.Range("I13").Value = 0.2184796247693259
$ws.Range("J13").Value = 0.2184796247693259
$ws.Range("M13").Value = 0.662052
$ws.Range("N13").Value = 1.986156
$ws.Range("O13").Value = 0.01300304006970129
$ws.Range("P13").Value = 0.0130030400697013
$ws.Range("Q13").Value = 12.162714198324
$ws.Range("R13").Value = 109.464427784916
$ws.Range("S13").Value = 0.002840899315288848
$ws.Range("T13").Value = 0.002840899315288849

$ws.Range("G14").Value = 25.652214
$ws.Range("H14").Value = 76.956642
$ws.Range("I14").Value = 0.3050685203844711
$ws.Range("J14").Value = 0.3050685203844711
$ws.Range("M14").Value = 46.63275166666667
$ws.Range("N14").Value = 139.898255
$ws.Range("O14").Value = 0.9158911059585902
$ws.Range("P14").Value = 0.9158911059585902
$ws.Range("Q14").Value = 1196.23332516219
$ws.Range("R14").Value = 10766.09992645971
$ws.Range("S14").Value = 0.279409544528084
$ws.Range("T14").Value = 0.279409544528084

$ws.Range("G15").Value = 25.652214
$ws.Range("H15").Value = 76.956642
$ws.Range("I15").Value = 0.3050685203844711
$ws.Range("J15").Value = 0.3050685203844711
$ws.Range("O15").Value = 0.05441917700612491
$ws.Range("P15").Value = 0.05441917700612491
$ws.Range("Q15").Value = 71.07617121633
$ws.Range("R15").Value = 639.6855409469699
$ws.Range("S15").Value = 0.01660157780979916
$ws.Range("T15").Value = 0.01660157780979916

$ws.Range("G16").Value = 25.652214
$ws.Range("H16").Value = 76.956642
$ws.Range("I16").Value = 0.3050685203844711
$ws.Range("J16").Value = 0.3050685203844711
$ws.Range("M16").Value = 0.8496050000000001
$ws.Range("N16").Value = 2.548815
$ws.Range("O16").Value = 0.01668667696558362
$ws.Range("P16").Value = 0.01668667696558362
$ws.Range("Q16").Value = 21.79424927547
$ws.Range("R16").Value = 196.14824347923
$ws.Range("S16").Value = 0.005090579852024233
$ws.Range("T16").Value = 0.005090579852024233

$ws.Range("G17").Value = 25.652214
$ws.Range("H17").Value = 76.956642
$ws.Range("I17").Value = 0.3050685203844711
$ws.Range("J17").Value = 0.3050685203844711
$ws.Range("M17").Value = 0.662052
$ws.Range("N17").Value = 1.986156
$ws.Range("O17").Value = 0.01300304006970129
$ws.Range("P17").Value = 0.0130030400697013
$ws.Range("Q17").Value = 16.983099583128
$ws.Range("R17").Value = 152.847896248152
$ws.Range("S17").Value = 0.003966818194563764
$ws.Range("T17").Value = 0.003966818194563765
